$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '26.125.27'
$ws.Range('E2').Value = '  +1.02%  '

$ws.Range('D3').Value = '1.605.86'
$ws.Range('E3').Value = '  +0.47%  '

$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.83%  '

$ws.Range('D5').Value = '212.78'
$ws.Range('E5').Value = '  +1.75%  '

$ws.Range('D6').Value = '1.00'
$ws.Range('E6').Value = '  -0.92%  '

$ws.Range('D7').Value = '0.482'
$ws.Range('E7').Value = '  +1.20%  '

$ws.Range('D8').Value = '0.250'
$ws.Range('E8').Value = '  +2.60%  '

$ws.Range('D9').Value = '0.0619'
$ws.Range('E9').Value = '  +1.36%  '

$ws.Range('D10').Value = '18.46'
$ws.Range('E10').Value = '  +3.60%  '

$ws.Range('D11').Value = '0.0790'
$ws.Range('E11').Value = '  +0.25%  '

$ws.Range('D12').Value = '1.830.07'
$ws.Range('E12').Value = '  +0.53%  '

$ws.Range('D13').Value = '1.593.69'
$ws.Range('E13').Value = '  -0.35%  '

$ws.Range('D14').Value = '4.06'
$ws.Range('E14').Value = '  +0.34%  '

$ws.Range('D15').Value = '0.513'
$ws.Range('E15').Value = '  +1.07%  '

$ws.Range('D16').Value = '26.149.77'
$ws.Range('E16').Value = '  +1.18%  '

$ws.Range('D17').Value = '0.0₃0731'
$ws.Range('E17').Value = '  +2.43%  '

$ws.Range('D18').Value = '60.80'
$ws.Range('E18').Value = '  +0.29%  '

$ws.Range('D19').Value = '0.999'
$ws.Range('E19').Value = '  -0.83%  '

$ws.Range('D20').Value = '198.17'
$ws.Range('E20').Value = '  +4.82%  '

$ws.Range('D21').Value = '4.27'
$ws.Range('E21').Value = '  +2.34%  '

$ws.Range('D22').Value = '9.50'
$ws.Range('E22').Value = '  +2.34%  '

$ws.Range('D23').Value = '6.01'
$ws.Range('E23').Value = '  +1.39%  '

$ws.Range('E24').Value = '  +2.70%  '

$ws.Range('D25').Value = '142.32'
$ws.Range('E25').Value = '  +0.50%  '

$ws.Range('D26').Value = '1.72'
$ws.Range('E26').Value = '  +1.51%  '

$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.81%  '

$ws.Range('D28').Value = '15.27'
$ws.Range('E28').Value = '  +2.46%  '

$ws.Range('D29').Value = '6.52'
$ws.Range('E29').Value = '  +0.32%  '

$ws.Range('E30').Value = '  -2.08%  '

$ws.Range('D31').Value = '0.0477'
$ws.Range('E31').Value = '  +2.01%  '

$ws.Range('D32').Value = '3.16'
$ws.Range('E32').Value = '  +2.76%  '

$ws.Range('D33').Value = '3.05'
$ws.Range('E33').Value = '  +2.01%  '

$ws.Range('D34').Value = '1.52'
$ws.Range('E34').Value = '  +4.38%  '

$ws.Range('E35').Value = '  -2.90%  '

$ws.Range('D36').Value = '1.106.71'
$ws.Range('E36').Value = '  -0.05%  '

$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').Value = '2.36'
$ws.Range('E37').Value = '  +0.49%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.0153'
$ws.Range('E38').Value = '  +2.19%  '

$ws.Range('D39').Value = '0.509'
$ws.Range('E39').Value = '  +3.10%  '

$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -1.04%  '

$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '0.791'
$ws.Range('E41').Value = '  -0.52%  '

$ws.Range('D42').Value = '0.802'
$ws.Range('E42').Value = '  +8.74%  '

$ws.Range('E43').Value = '  +1.86%  '

$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.734.60'
$ws.Range('E44').Value = '  +0.06%  '

$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '93.03'
$ws.Range('E45').Value = '  -2.64%  '

$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₆0107'
$ws.Range('E46').Value = '  -5.34%  '

$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = '1.56'
$ws.Range('E47').Value = '  +7.48%  '

$ws.Range('D48').Value = '53.95'
$ws.Range('E48').Value = '  +1.87%  '

$ws.Range('D49').Value = '0.0508'
$ws.Range('E49').Value = '  -0.28%  '

$ws.Range('D50').Value = '0.406'
$ws.Range('E50').Value = '  -0.83%  '

$ws.Range('D51').Value = '0.999'
$ws.Range('E51').Value = '  -0.87%  '
